# Generate Report for Handoff
#
# Marks the e53aacaa-3cc5-4a25-a9ee-0e9c9e7ee0ef.md source file as handed off
# ("Ready for handoff") across the Overview, zh-cn and de-de sheets, updating
# the relevant priority / handoff-datetime fields to match.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-19 14:13:49"

# The longer "Ready for handoff" status text makes Excel widen the
# zh-cn/de-de status columns (autofit); reproduce the resulting width.
$overview.Columns.Item(5).ColumnWidth = 16.38
$overview.Columns.Item(6).ColumnWidth = 16.38

# ---- zh-cn sheet -------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("E3").Value = "mt"
$zhcn.Range("H3").Value = "2016-08-19 14:13:45"
$zhcn.Columns.Item(3).ColumnWidth = 16.38

# ---- de-de sheet -------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("E3").Value = "mt"
$dede.Range("H3").Value = "2016-08-19 14:13:49"
$dede.Columns.Item(3).ColumnWidth = 16.38
